$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.328.97"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "1.869.81"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("D5").Value = "'235.07"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.4706"
$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("D8").Value = "'0.2870"
$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").Value = "'0.06571"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").Value = "'21.61"
$ws.Range("E10").Value = "  -3.22%  "

$ws.Range("D11").Value = "'0.08016"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("E12").Value = "  -0.84%  "

$ws.Range("D13").Value = "1.869.69"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("E14").Value = "  -1.31%  "

$ws.Range("D15").Value = "'0.6839"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").Value = "'269.42"

$ws.Range("D17").Value = "30.319.68"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("E18").Value = "  +2.20%  "

$ws.Range("D19").Value = "'0.000007626"
$ws.Range("E19").Value = "  +3.90%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "2.113.54"
$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").Value = "'5.282"
$ws.Range("E23").Value = "  -2.25%  "

$ws.Range("D24").Value = "'6.216"
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("D25").Value = "'9.427"
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("D26").Value = "'168.84"
$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("D27").Value = "'18.91"
$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("D29").Value = "'1.370"
$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("D30").Value = "'0.09931"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").Value = "'4.364"
$ws.Range("E31").Value = "  -0.82%  "

$ws.Range("D32").Value = "'1.461"
$ws.Range("E32").Value = "  -1.37%  "

$ws.Range("D33").Value = "'4.074"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").Value = "'0.04702"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").Value = "'0.7008"
$ws.Range("E36").Value = "  -1.63%  "

$ws.Range("D37").Value = "'2.702"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").Value = "'0.01874"
$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("D39").Value = "'2.636"
$ws.Range("E39").Value = "  +0.83%  "

$ws.Range("D40").Value = "'6.301"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = "'71.84"
$ws.Range("E41").Value = "  -6.67%  "

$ws.Range("D42").Value = "'1.967"
$ws.Range("E42").Value = "  +0.29%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8415"
$ws.Range("E43").Value = "  -1.14%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4167"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").Value = "'102.84"
$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("D47").Value = "'9.208"
$ws.Range("E47").Value = "  -1.47%  "

$ws.Range("D48").Value = "'7.051"

$ws.Range("D49").Value = "'913.71"
$ws.Range("E49").Value = "  -5.47%  "

$ws.Range("D50").Value = "'34.48"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").Value = "'0.05698"
$ws.Range("E51").Value = "  +0.91%  "
